# Burndown chart update: refresh the Sprint 2 and Release burndown figures,
# remove the stray "ALSO DO A SPRINT 2 BURNDOWN AS WELL" note, and move the
# selection/view to where work left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burn down chart")

# --- Release Burndown (rows 3-5) ---
$ws.Range("B3").Value = 31
$ws.Range("C3").Value = 31
$ws.Range("C4").Value = 19

# --- Sprint 2 Burndown (rows 13-20) ---
$ws.Range("B13").Value = 77
$ws.Range("C13").Value = 77
$ws.Range("C14").Value = 75
$ws.Range("C15").Value = 71
$ws.Range("C16").Value = 70
$ws.Range("C17").Value = 67
$ws.Range("C18").Value = 60
$ws.Range("C19").Value = 52
$ws.Range("C20").Value = 48

# Remove the leftover reminder note - it's done, so delete the cell entirely.
$ws.Range("A36").Clear()

# Leave the selection where work continues.
$ws.Range("A36").Select()
